# Update division-problem answers throughout the document.
# Each original answer string is unique in the document, so a plain
# Find/Replace (no wildcards) for each pair is sufficient and safe.

$d = $word.ActiveDocument

$replacements = @(
    @{ old = "456÷4=114, 0"; new = "552÷5=110, 2" },
    @{ old = "987÷5=197, 2"; new = "567÷4=141, 3" },
    @{ old = "999÷6=166, 3"; new = "403÷6=67, 1" },
    @{ old = "450÷2=225, 0"; new = "649÷8=81, 1" },
    @{ old = "634÷4=158, 2"; new = "150÷4=37, 2" },
    @{ old = "246÷5=49, 1";  new = "810÷7=115, 5" },
    @{ old = "142÷8=17, 6";  new = "166÷7=23, 5" },
    @{ old = "915÷5=183, 0"; new = "487÷5=97, 2" },
    @{ old = "460÷4=115, 0"; new = "809÷7=115, 4" },
    @{ old = "457÷5=91, 2";  new = "602÷2=301, 0" },
    @{ old = "199÷6=33, 1";  new = "551÷6=91, 5" },
    @{ old = "977÷2=488, 1"; new = "352÷9=39, 1" },
    @{ old = "493÷5=98, 3";  new = "927÷6=154, 3" },
    @{ old = "703÷3=234, 1"; new = "209÷3=69, 2" },
    @{ old = "979÷5=195, 4"; new = "186÷2=93, 0" },
    @{ old = "622÷2=311, 0"; new = "743÷3=247, 2" },
    @{ old = "158÷7=22, 4";  new = "881÷5=176, 1" },
    @{ old = "438÷4=109, 2"; new = "186÷4=46, 2" },
    @{ old = "797÷8=99, 5";  new = "892÷9=99, 1" },
    @{ old = "565÷6=94, 1";  new = "351÷4=87, 3" },
    @{ old = "775÷4=193, 3"; new = "898÷2=449, 0" },
    @{ old = "776÷6=129, 2"; new = "600÷3=200, 0" },
    @{ old = "569÷2=284, 1"; new = "118÷5=23, 3" },
    @{ old = "723÷8=90, 3";  new = "271÷8=33, 7" },
    @{ old = "778÷5=155, 3"; new = "213÷8=26, 5" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
